$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Append new row 18 to the Logs sheet (new test-mail entry)
$logs.Range("A18").Value = "Kun je nagaan of dit nog leverbaar is?"
$logs.Range("B18").Value = "mailmind.test@zohomail.eu"
$logs.Range("C18").Value = "Testmail #8: Kun je nagaan of dit nog leverbaar is?"
$logs.Range("D18").Value = "Inkoop / Bestellingen"
$logs.Range("E18").Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Range("F18").Value = "2025-08-03 14:49:01"
$logs.Range("G18").Value = "Ja"
$logs.Range("H18").Value = "Ja"
$logs.Range("I18").Value = "Nee"
$logs.Range("J18").Value = "Nee"

# Extend the existing conditional-formatting rules so they keep covering the
# whole data range now that row 18 was added.
$logs.Range("D2:D17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D18"))
$logs.Range("G2:G17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G18"))
$logs.Range("H2:H17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H18"))
$logs.Range("I2:I17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I18"))
$logs.Range("J2:J17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J18"))

# Update the Dashboard's count for "Inkoop / Bestellingen" (row 5) to reflect
# the newly added log entry.
$dash.Range("B5").Value = 3
